$wb = $excel.ActiveWorkbook

# ---------- Sheet: ALC ----------
$ws = $wb.Worksheets.Item("ALC")

# row 21
$ws.Range("H21").Value = 19499.5
$ws.Range("J21").Value = 19499.5
$ws.Range("L21").Value = 19499.5
$ws.Range("N21").Value = -20435.5

# row 23
$ws.Range("H23").Value = 19499.5
$ws.Range("J23").Value = 19499.5
$ws.Range("L23").Value = 19499.5
$ws.Range("N23").Value = -19967.5

# row 34 (N34 removed, M34 updated)
$ws.Range("H34").Value = 975
$ws.Range("I34").Value = 975
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 975
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -772
$ws.Range("N34").ClearContents()

# row 36 (N36 removed, M36 updated)
$ws.Range("H36").Value = 975
$ws.Range("I36").Value = 975
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 975
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -260
$ws.Range("N36").ClearContents()

# row 38
$ws.Range("H38").Value = 1428629.9
$ws.Range("I38").Value = 68.166664
$ws.Range("J38").Value = 10000000
$ws.Range("K38").Value = 204.499992
$ws.Range("L38").Value = 30000000
$ws.Range("M38").Value = 167.500008
$ws.Range("N38").Value = -30000744

# row 113
$ws.Range("H113").Value = 2800.7896
$ws.Range("I113").Value = 2473.5
$ws.Range("K113").Value = 2473.5
$ws.Range("M113").Value = 780.5

# row 132
$ws.Range("H132").Value = 3380615.8
$ws.Range("I132").Value = 2051.652
$ws.Range("J132").Value = 50004800
$ws.Range("K132").Value = 6154.956
$ws.Range("L132").Value = 150014400
$ws.Range("M132").Value = -3624.956
$ws.Range("N132").Value = -150019460

# row 138
$ws.Range("H138").Value = 5294589
$ws.Range("I138").Value = 23813310
$ws.Range("J138").Value = 3526.0613
$ws.Range("K138").Value = 71439930
$ws.Range("L138").Value = 10578.1839
$ws.Range("M138").Value = -71434790
$ws.Range("N138").Value = -20858.1839

# ---------- Sheet: ARM ----------
$ws = $wb.Worksheets.Item("ARM")

# row 32
$ws.Range("H32").Value = 7200.79
$ws.Range("I32").Value = 5537.0933
$ws.Range("J32").Value = 17420.643
$ws.Range("K32").Value = 5537.0933
$ws.Range("L32").Value = 17420.643
$ws.Range("M32").Value = -5250.0933
$ws.Range("N32").Value = -17994.643

# row 33 (new M33 added)
$ws.Range("H33").Value = 19800
$ws.Range("I33").Value = 19800
$ws.Range("K33").Value = 19800
$ws.Range("M33").Value = -19471

# row 37
$ws.Range("H37").Value = 8729.777
$ws.Range("J37").Value = 8729.777
$ws.Range("L37").Value = 8729.777
$ws.Range("N37").Value = -9275.777

# row 43
$ws.Range("H43").Value = 14310.8
$ws.Range("J43").Value = 14310.8
$ws.Range("L43").Value = 14310.8
$ws.Range("N43").Value = -14936.8

# row 45
$ws.Range("H45").Value = 2166043.2
$ws.Range("I45").Value = 2842510
$ws.Range("J45").Value = 1350
$ws.Range("K45").Value = 2842510
$ws.Range("L45").Value = 1350
$ws.Range("M45").Value = -2842133
$ws.Range("N45").Value = -2104

# row 55 (M55 removed)
$ws.Range("H55").Value = 28788.273
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 28788.273
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 28788.273
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -29418.273

# row 110
$ws.Range("H110").Value = 1629.4286
$ws.Range("I110").Value = 1240
$ws.Range("J110").Value = 2603
$ws.Range("K110").Value = 1240
$ws.Range("L110").Value = 2603
$ws.Range("M110").Value = 805
$ws.Range("N110").Value = -6693

# row 122
$ws.Range("H122").Value = 2736.6667
$ws.Range("I122").Value = 2736.9333
$ws.Range("J122").Value = 2735.3333
$ws.Range("K122").Value = 8210.7999
$ws.Range("L122").Value = 8205.999899999999
$ws.Range("M122").Value = -5760.7999
$ws.Range("N122").Value = -13105.9999

# ---------- Sheet: BSM ----------
$ws = $wb.Worksheets.Item("BSM")

# row 56
$ws.Range("H56").Value = 50110
$ws.Range("J56").Value = 50110
$ws.Range("L56").Value = 50110
$ws.Range("N56").Value = -51588

# row 134
$ws.Range("H134").Value = 3773.4424
$ws.Range("I134").Value = 2605.258
$ws.Range("J134").Value = 5497.905
$ws.Range("K134").Value = 7815.773999999999
$ws.Range("L134").Value = 16493.715
$ws.Range("M134").Value = -5280.773999999999
$ws.Range("N134").Value = -21563.715

# ---------- Sheet: CRP ----------
$ws = $wb.Worksheets.Item("CRP")

# row 54
$ws.Range("H54").Value = 30055.2
$ws.Range("J54").Value = 30055.2
$ws.Range("L54").Value = 30055.2
$ws.Range("N54").Value = -31371.2

# row 141
$ws.Range("H141").Value = 37018.46
$ws.Range("J141").Value = 37586.086
$ws.Range("L141").Value = 37586.086
$ws.Range("N141").Value = -47946.086

# ---------- Sheet: CUL ----------
$ws = $wb.Worksheets.Item("CUL")

# row 13 (N13 removed, M13 updated)
$ws.Range("H13").Value = 420
$ws.Range("I13").Value = 420
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1260
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1092
$ws.Range("N13").ClearContents()

# row 102
$ws.Range("H102").Value = 10800
$ws.Range("J102").Value = 10800
$ws.Range("L102").Value = 32400
$ws.Range("N102").Value = -37268

# row 107
$ws.Range("H107").Value = 756.0789
$ws.Range("I107").Value = 229
$ws.Range("J107").Value = 875.0968
$ws.Range("K107").Value = 687
$ws.Range("L107").Value = 2625.2904
$ws.Range("M107").Value = 1233
$ws.Range("N107").Value = -6465.2904

# row 129
$ws.Range("H129").Value = 3161.739
$ws.Range("I129").Value = 1102.8572
$ws.Range("J129").Value = 4062.5
$ws.Range("K129").Value = 3308.5716
$ws.Range("L129").Value = 12187.5
$ws.Range("M129").Value = 1691.4284
$ws.Range("N129").Value = -22187.5

# row 131
$ws.Range("H131").Value = 2949.5254
$ws.Range("I131").Value = 804.2143
$ws.Range("J131").Value = 3616.9556
$ws.Range("K131").Value = 2412.6429
$ws.Range("L131").Value = 10850.8668
$ws.Range("M131").Value = 2627.3571
$ws.Range("N131").Value = -20930.8668

# row 137
$ws.Range("H137").Value = 3011.465
$ws.Range("I137").Value = 2790.6667
$ws.Range("J137").Value = 3129.75
$ws.Range("K137").Value = 8372.000100000001
$ws.Range("L137").Value = 9389.25
$ws.Range("M137").Value = -3272.000100000001
$ws.Range("N137").Value = -19589.25

# ---------- Sheet: GSM ----------
$ws = $wb.Worksheets.Item("GSM")

# row 41 (new N41 added)
$ws.Range("H41").Value = 1066.6666
$ws.Range("I41").Value = 350
$ws.Range("J41").Value = 2500
$ws.Range("K41").Value = 350
$ws.Range("L41").Value = 2500
$ws.Range("M41").Value = 5
$ws.Range("N41").Value = -3210

# ---------- Sheet: LTW ----------
$ws = $wb.Worksheets.Item("LTW")

# row 51
$ws.Range("H51").Value = 22600
$ws.Range("J51").Value = 22600
$ws.Range("L51").Value = 22600
$ws.Range("N51").Value = -23556

# row 122
$ws.Range("H122").Value = 7449.8335
$ws.Range("I122").Value = 3479
$ws.Range("J122").Value = 11420.667
$ws.Range("K122").Value = 10437
$ws.Range("L122").Value = 34262.001
$ws.Range("M122").Value = -7987
$ws.Range("N122").Value = -39162.001

# ---------- Sheet: WVR ----------
$ws = $wb.Worksheets.Item("WVR")

# row 122
$ws.Range("H122").Value = 3025.08
$ws.Range("I122").Value = 2839.7144
$ws.Range("J122").Value = 3998.25
$ws.Range("K122").Value = 8519.143199999999
$ws.Range("L122").Value = 11994.75
$ws.Range("M122").Value = -6069.143199999999
$ws.Range("N122").Value = -16894.75

# row 140
$ws.Range("H140").Value = 111666.664
$ws.Range("J140").Value = 111666.664
$ws.Range("L140").Value = 111666.664
$ws.Range("N140").Value = -122026.664
